$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("M2").Value = 4.858785666666667
$ws.Range("N2").Value = 14.576357
$ws.Range("O2").Value = 0.383527542896884
$ws.Range("P2").Value = 0.383527542896884
$ws.Range("Q2").Value = 0.331101949255
$ws.Range("R2").Value = 2.979917543295
$ws.Range("S2").Value = 0.383527542896884
$ws.Range("T2").Value = 0.383527542896884

# Row 3
$ws.Range("N3").Value = 7.354217
$ws.Range("O3").Value = 0.193501351259474
$ws.Range("P3").Value = 0.193501351259474
$ws.Range("S3").Value = 0.193501351259474
$ws.Range("T3").Value = 0.193501351259474

# Row 4
$ws.Range("M4").Value = 0.9922233333333333
$ws.Range("N4").Value = 2.97667
$ws.Range("O4").Value = 0.07832100511224221
$ws.Range("P4").Value = 0.07832100511224221
$ws.Range("Q4").Value = 0.06761505905
$ws.Range("R4").Value = 0.6085355314500001
$ws.Range("S4").Value = 0.07832100511224221
$ws.Range("T4").Value = 0.07832100511224221

# Row 5
$ws.Range("M5").Value = 0.8999579999999999
$ws.Range("N5").Value = 2.699874
$ws.Range("O5").Value = 0.07103805438843064
$ws.Range("P5").Value = 0.07103805438843064
$ws.Range("Q5").Value = 0.06132763790999999
$ws.Range("R5").Value = 0.55194874119
$ws.Range("S5").Value = 0.07103805438843064
$ws.Range("T5").Value = 0.07103805438843064

# Row 6
$ws.Range("M6").Value = 2.962941666666666
$ws.Range("N6").Value = 8.888824999999999
$ws.Range("O6").Value = 0.2338793713333444
$ws.Range("P6").Value = 0.2338793713333444
$ws.Range("Q6").Value = 0.201909659875
$ws.Range("R6").Value = 1.817186938875
$ws.Range("S6").Value = 0.2338793713333444
$ws.Range("T6").Value = 0.2338793713333444

# Row 7
$ws.Range("M7").Value = 0.5033603333333333
$ws.Range("N7").Value = 1.510081
$ws.Range("O7").Value = 0.03973267500962479
$ws.Range("P7").Value = 0.0397326750096248
$ws.Range("Q7").Value = 0.034301489915
$ws.Range("R7").Value = 0.308713409235
$ws.Range("S7").Value = 0.03973267500962479
$ws.Range("T7").Value = 0.0397326750096248
